# Update "想去人数" (want-to-go count) figures in the F column.
# Sheet "展览" (Exhibition) rows 5,7-12
# Sheet "全部类型" (All types) rows 6,8,11-15

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 243
$wsExhibition.Range("F7").Value = 235
$wsExhibition.Range("F8").Value = 2268
$wsExhibition.Range("F9").Value = 383
$wsExhibition.Range("F10").Value = 5615
$wsExhibition.Range("F11").Value = 133
$wsExhibition.Range("F12").Value = 367

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 243
$wsAll.Range("F8").Value = 235
$wsAll.Range("F11").Value = 2268
$wsAll.Range("F12").Value = 383
$wsAll.Range("F13").Value = 5615
$wsAll.Range("F14").Value = 133
$wsAll.Range("F15").Value = 367
